# "Generate Report for Handback"
#
# The localization status workbook gets a handback-transform report:
#   - the "Ready for handoff" status cells flip to "Handback transform failed"
#     on the Overview sheet (cols E/F) and on each language sheet (col C)
#   - the (previously empty) Error Detail column (col P) on each language
#     sheet is filled in with the specific handback/handoff filename mismatch
#   - a couple of columns are widened so the new, longer text is readable

$wb = $excel.ActiveWorkbook

$failedStatus = "Handback transform failed"

$zhError = "Handback file name: 41tcgxug.j4t is different with handoff file name: 9c0d7686-421c-4e84-90ef-f0810e5d5884.3a385625cbe39d4ac36af19f237f7253341ed15e.zh-cn."
$deError = "Handback file name: 41tcgxug.j4t is different with handoff file name: 9c0d7686-421c-4e84-90ef-f0810e5d5884.3a385625cbe39d4ac36af19f237f7253341ed15e.de-de."

# --- Overview sheet: status columns E2/F2 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $failedStatus
$wsOverview.Range("F2").Value = $failedStatus

# Widen columns E and F to fit the longer status text (COM ColumnWidth
# snaps to the nearest 1/6-character pixel grid, so feed the value that
# lands on the grid point closest to the target stored width).
$wsOverview.Columns.Item(5).ColumnWidth = 23.833333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 23.833333333333336

# --- zh-cn sheet: status column C2, error detail column P2 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $failedStatus
$wsZhCn.Range("P2").Value = $zhError

$wsZhCn.Columns.Item(3).ColumnWidth = 23.833333333333336
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: status column C2, error detail column P2 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $failedStatus
$wsDeDe.Range("P2").Value = $deError

$wsDeDe.Columns.Item(3).ColumnWidth = 23.833333333333336
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
